$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: Hydrogen -> clear D3 (was 736.9689587095404), leave as an empty cell
$ws.Range("D3").Value = ""
$ws.Range("D3").Style = "Normal"

# --- Row 4: Methanol -> C4 changes from 17437.30652341867 to 0
$ws.Range("C4").Value = 0

# --- Row 5: Ammonia -> C5 changes from 71177.51806802199 to 1237.087816454243
$ws.Range("C5").Value = 1237.087816454243

# --- Row 6: Biomass -> unchanged

# --- Row 7: relabel "Other" -> "Biogas", update D7 value
$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value = 945.213535710629

# --- Row 8 (new): "Other" row, mirroring the style/layout of row 7
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A8").Value = "Other"

$ws.Range("B8").Value = ""
$ws.Range("B8").Style = "Normal"

$ws.Range("C8").Value = ""
$ws.Range("C8").Style = "Normal"

$ws.Range("D8").Value = 732.0507850773076
